$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.176.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").Value = "'1.871.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'307.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.5066"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("D8").Value = "'0.3763"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").Value = "'0.07163"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("D10").Value = "'0.8910"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").Value = "'0.07583"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "'5.334"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("D15").Value = "'89.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "'0.000008542"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").Value = "'14.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.88%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'27.222.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("E22").Value = "  -1.90%  "
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("D24").Value = "'6.504"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("D25").Value = "'150.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("D26").Value = "'1.845"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("D27").Value = "'18.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "'2.115"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.70%  "
$ws.Range("D29").Value = "'112.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.00%  "
$ws.Range("E30").Value = "  -2.87%  "
$ws.Range("D31").Value = "'4.727"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("D32").Value = "'0.08994"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D34").Value = "'3.095"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.53%  "
$ws.Range("D35").Value = "'0.7522"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("E36").Value = "  -4.12%  "
$ws.Range("E37").Value = "  +1.63%  "
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("D40").Value = "'1.075"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("E41").Value = "  -2.67%  "
$ws.Range("D42").Value = "'6.632"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("D43").Value = "'115.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.61%  "
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("D46").Value = "'0.4667"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.64%  "
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("E48").Value = "  -4.12%  "
$ws.Range("D49").Value = "'1.571"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("D50").Value = "'65.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.25%  "
$ws.Range("E51").Value = "  -0.32%  "
